# 5S Report Template: add a number format to the Score header/data column,
# fill the (previously empty) Score cells with a single space, shrink the
# table from 8 data rows to 3, and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Score" header cell (E2) now shows a 0.00 number format
$ws.Range("E2").NumberFormat = "0.00"

# The three remaining score cells get a single-space placeholder value
$ws.Range("E3").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("E5").Value = " "

# Remove the trailing empty rows 6-8 (table now only spans to row 5)
$ws.Rows("6:8").Delete()

# Move the active selection to H9
$ws.Range("H9").Select() | Out-Null
